$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

$ws.Range("A8").Value = "Beschadigd product ontvangen"
$ws.Range("B8").Value = "mailmind.test@zohomail.eu"
$ws.Range("C8").Value = "Het product dat ik heb ontvangen is beschadigd aangekomen."
$ws.Range("D8").Value = "Retour / Terugbetaling"
$ws.Range("E8").Value = "Beste klant,`nBedankt voor het melden van de beschadiging van het ontvangen product. Om dit probleem snel voor u op te lossen, heb ik wat extra informatie nodig. Kunt u alstublieft een foto van de beschadiging van het product meesturen? Dit zal ons helpen om de situatie beter te begrijpen en een passende oplossing te bieden.`nIk kijk uit naar uw antwoord.`nMet vriendelijke groet,`n[Naam] `nE-mailassistent"
$ws.Range("F8").Value = "2025-06-24 19:48:56"
$ws.Range("G8").Value = "Ja"
$ws.Rows.Item(8).AutoFit()

$dash.Range("A2").Value = "Retour / Terugbetaling"
$dash.Range("B2").Value = 2
$dash.Range("A3").Value = "Factuur / Administratie"
$dash.Range("B3").Value = 2

$ws.Range("D2:D7").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D8"))
$ws.Range("G2:G7").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G8"))
